# Applies the "Add files via upload" edit to SecuritiesEnforcement.xlsx:
# fills in the Outcome..SEC Office columns (D:M) for four enforcement-action
# rows (34-37) that previously only had Date/Case/Description (A:C) filled
# in, fixes up a stray missing style on G31, and updates the sheet's
# scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: G31 was missing the shared Times-New-Roman cell style that the
# rest of the row uses (s="2"); pick it up from a neighboring cell so the
# same style index is reused instead of a brand-new one being minted.
$ws.Range("F31").Copy() | Out-Null
$ws.Range("G31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 34: Salt Blockchain Inc. settlement ---------------------------
# (H before G to match the shared-string insertion order in the workbook)
$ws.Range("H34").Value = "Salt Blockchain, Inc."
$ws.Range("G34").Value = "SALT"
$ws.Range("D34").Value = "Settlement"
$ws.Range("E34").Value = "Unregistered Offering"
$ws.Range("F34").Value = "Civil"
$ws.Range("I34").Value = "Ethereum"
$ws.Range("J34").Value = 47000000
$ws.Range("K34").Value = 1
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = "Denver"

# --- Row 35: SoluTech, Inc., et al. settlement --------------------------
$ws.Range("D35").Value = "Settlement"
$ws.Range("E35").Value = "Unregistered Offering"
$ws.Range("F35").Value = "Civil"
$ws.Range("G35").Value = "SCRL"
$ws.Range("H35").Value = "SoluTech, Inc., et al."
$ws.Range("I35").Value = "Ethereum"
$ws.Range("J35").Value = 2400000
$ws.Range("K35").Value = 1
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = "Washington, D.C."

# --- Row 36: Unikrn, Inc. settlement -------------------------------------
$ws.Range("D36").Value = "Settlement"
$ws.Range("E36").Value = "Unregistered Offering"
$ws.Range("F36").Value = "Civil"
$ws.Range("G36").Value = "UKG"
$ws.Range("H36").Value = "Unikrn, Inc."
$ws.Range("I36").Value = "Ethereum"
$ws.Range("J36").Value = 6100000
$ws.Range("K36").Value = 1
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = "Washington"

# --- Row 37: SEC v. FLiK, et al. -----------------------------------------
$ws.Range("D37").Value = "Settlement"
$ws.Range("E37").Value = "Unregistered Offering"
$ws.Range("F37").Value = "Civil and Criminal"
$ws.Range("H37").Value = "FLiK and CoinSpark"
$ws.Range("I37").Value = "Ethereum"
$ws.Range("J37").Value = 2200000
$ws.Range("K37").Value = 1
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = "Northern Georgia"
$ws.Range("G37").Value = "FLiK"

# --- Sheet view: drop the scrolled "G1" top-left cell and move the
# active selection to P34 -------------------------------------------------
$ws.Range("P34").Select() | Out-Null
